# Updated symbol list on Tue Dec 27 20:06:36 UTC 2022 with GitHub Actions
#
# The sheet stores every data cell as text (inline strings), including
# values that look numeric (price, hour). Plain `Range.Value = "20"`
# assignments would be auto-coerced to a Number by Excel, which would
# change the cell's stored type from Text -> Number. To keep them as
# literal text (matching the source data) we prefix numeric-looking
# values with a leading apostrophe, exactly like typing '20 into Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $value)
    $ws.Range($addr).Value = "'" + $value
}

# Row 2 - BNB: price + hour refresh
Set-TextValue "D2" "245.96"
Set-TextValue "G2" "20"

# Row 3 - OKB: hour refresh only
Set-TextValue "G3" "20"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.371"
Set-TextValue "G4" "20"

# Row 5 - Cronos
Set-TextValue "D5" "0.05836"
Set-TextValue "G5" "20"

# Row 6 - KuCoinToken
Set-TextValue "D6" "6.469"
Set-TextValue "G6" "20"

# Row 7 - GateToken
Set-TextValue "D7" "3.354"
Set-TextValue "G7" "20"

# Row 8 - MXToken
Set-TextValue "D8" "0.8103"
Set-TextValue "G8" "20"

# Row 9 - FTXToken
Set-TextValue "D9" "0.9202"
Set-TextValue "G9" "20"

# Row 10 - WazirX
Set-TextValue "D10" "0.1410"
Set-TextValue "G10" "20"

# Row 11 - MandalaExchangeToken
Set-TextValue "D11" "0.07378"
Set-TextValue "G11" "20"

# Row 12 - LiechtensteinCryptoassetsExchange: hour refresh only
Set-TextValue "G12" "20"

# Row 13 - BitrueCoin
Set-TextValue "D13" "0.03026"
Set-TextValue "G13" "20"

# Row 14 - BitMartToken
Set-TextValue "D14" "0.09374"
Set-TextValue "G14" "20"

# Row 15 - MCDex
Set-TextValue "D15" "3.848"
Set-TextValue "G15" "20"

# Row 16 - BitForexToken
Set-TextValue "D16" "0.001558"
Set-TextValue "G16" "20"

# Row 17 - CoinExToken
Set-TextValue "D17" "0.04691"
Set-TextValue "G17" "20"

# Row 18 - One: price + volume label lost its "Worstin24h" badge + hour refresh
Set-TextValue "D18" "0.0005978"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "G18" "20"

# Row 19 - TigerCash
Set-TextValue "D19" "0.006008"
Set-TextValue "G19" "20"

# Row 20 - BitKan: hour refresh only
Set-TextValue "G20" "20"

# Row 21 - HotbitToken
Set-TextValue "D21" "0.004689"
Set-TextValue "G21" "20"

# Row 22 - NitroEx
Set-TextValue "D22" "0.00008800"
Set-TextValue "G22" "20"

# Row 23 - LEO
Set-TextValue "D23" "3.593"
Set-TextValue "G23" "20"

# Rows 24-27 - BTSEToken, BitpandaEcosystemToken, ProBitToken, AAXToken: hour refresh only
Set-TextValue "G24" "20"
Set-TextValue "G25" "20"
Set-TextValue "G26" "20"
Set-TextValue "G27" "20"

# Row 28 - UpBots
Set-TextValue "D28" "0.0002349"
Set-TextValue "G28" "20"

# Rows 29-39 - hour refresh only (prices remain "--")
Set-TextValue "G29" "20"
Set-TextValue "G30" "20"
Set-TextValue "G31" "20"
Set-TextValue "G32" "20"
Set-TextValue "G33" "20"
Set-TextValue "G34" "20"
Set-TextValue "G35" "20"
Set-TextValue "G36" "20"
Set-TextValue "G37" "20"
Set-TextValue "G38" "20"
Set-TextValue "G39" "20"

# Row 40 - IDEX
Set-TextValue "D40" "0.03841"
Set-TextValue "G40" "20"

# Rows 41-43 reshuffled: KickToken/BKEXToken/CEJI -> BKEXToken/CEJI/KickToken
# Row 41 now holds BKEXToken's data (rank index 39 stays put in column A)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1066"
$ws.Range("E41").Value = "40BKEXTokenBKK"
Set-TextValue "G41" "20"

# Row 42 now holds CEJI's data
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002750"
$ws.Range("E42").Value = "41CEJICEJI"
Set-TextValue "G42" "20"

# Row 43 now holds KickToken's data (and gains the "Worstin24h" badge)
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003086"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
Set-TextValue "G43" "20"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008362"
Set-TextValue "G44" "20"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005261"
Set-TextValue "G45" "20"

# Row 46 - Kangarootoken: hour refresh only
Set-TextValue "G46" "20"

# Row 47 - CoinbaseStockToken
Set-TextValue "D47" "0.7097"
Set-TextValue "G47" "20"

# Row 48 - BOLO
Set-TextValue "D48" "0.001837"
Set-TextValue "G48" "20"

# Row 49 - CryptobidCoin
Set-TextValue "D49" "0.00002099"
Set-TextValue "G49" "20"

# Row 50 - SpecialPowerGold
Set-TextValue "D50" "0.0001999"
Set-TextValue "G50" "20"

# Row 51 - DigiFinexToken: hour refresh only
Set-TextValue "G51" "20"
